$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell text updates (row 1) ---
# Original: C1="Enseignant", D1="Nombre d'heures"
# Target:   C1="Chef  Module" (note: two spaces), D1="Composants"
# Set D1 first, then C1, so the shared-string table ends up with
# "Composants" at index 2 and "Chef  Module" at index 3 (matching
# the target sharedStrings.xml / cell <v> index layout).
$ws.Range("D1").Value = "Composants"
$ws.Range("C1").Value = "Chef  Module"

# --- Column widths ---
# The engine quantizes the stored column width to steps of 1/6 of a
# character, so these inputs are chosen to land safely in the middle
# of the bucket that rounds to the target stored widths (35 exactly
# for column C; the closest representable stored width to the target
# 24.5703125 for column D is 24.5).
$ws.Columns.Item(3).ColumnWidth = 34.16
$ws.Columns.Item(4).ColumnWidth = 23.665

# --- Selection change ---
$ws.Range("E8").Select() | Out-Null
